$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new entry (row 5) to the "Learning goal 1" hour registration table:
# Date 17 Dec 2021 (serial 44547), From-to 14.00 - 15.00, Hours 1h,
# Activity "Working on movement".
# Shared strings must be created in this order so new <si> entries land as
# "Working on movement" (index 18) then "14.00 - 15.00" (index 19).

$ws.Range("C5").Value = 44547
$ws.Range("C5").NumberFormat = "d-mmm"

$ws.Range("G5").Value = "Working on movement"
$ws.Range("D5").Value = "14.00 - 15.00"
$ws.Range("F5").Value = "1h"

$null = $ws.Range("I23").Select()
